# Re-apply the latest cryptos.xlsx price/volume snapshot (GitHub Actions refresh).
# D/E columns hold plain display text (prices use "." as a thousands separator in
# this sheet, e.g. "30.248.05", and volumes are "  +1.29%  " with padding spaces),
# so values that look numeric are entered with a leading apostrophe to force Excel
# to store them as text instead of auto-converting to a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.248.05'
$ws.Range("E2").Value = '  +1.29%  '

$ws.Range("D3").Value = '2.085.14'
$ws.Range("E3").Value = '  -1.42%  '

$ws.Range("D4").Value = '''1.006'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").Value = '''340.02'
$ws.Range("E5").Value = '  -2.25%  '

$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").Value = '''0.5273'
$ws.Range("E7").Value = '  +1.67%  '

$ws.Range("D8").Value = '''0.4382'
$ws.Range("E8").Value = '  -1.89%  '

$ws.Range("D9").Value = '''54.96'
$ws.Range("E9").Value = '  +1.40%  '

$ws.Range("D10").Value = '''0.09335'
$ws.Range("E10").Value = '  -0.34%  '

$ws.Range("D11").Value = '''1.174'
$ws.Range("E11").Value = '  -0.67%  '

$ws.Range("E12").Value = '  -2.74%  '

$ws.Range("D13").Value = '''8.477'
$ws.Range("E13").Value = '  +0.98%  '

# Row 14: coin swapped position with the adjacent row
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '2.091.32'
$ws.Range("E14").Value = '  -0.71%  '

# Row 15: coin swapped position with the adjacent row
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '''6.854'
$ws.Range("E15").Value = '  +0.12%  '

$ws.Range("D16").Value = '''101.56'

$ws.Range("D17").Value = '''0.00001156'
$ws.Range("E17").Value = '  -0.96%  '

$ws.Range("E18").Value = '  -0.16%  '

$ws.Range("D19").Value = '''21.01'
$ws.Range("E19").Value = '  -2.69%  '

$ws.Range("D20").Value = '''0.06695'
$ws.Range("E20").Value = '  +0.35%  '

$ws.Range("D21").Value = '''6.290'
$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("E22").Value = '  -0.24%  '

$ws.Range("D23").Value = '30.263.43'
$ws.Range("E23").Value = '  +1.18%  '

$ws.Range("D24").Value = '''12.38'
$ws.Range("E24").Value = '  -2.57%  '

$ws.Range("D25").Value = '''2.323'
$ws.Range("E25").Value = '  -0.24%  '

$ws.Range("D26").Value = '''21.75'
$ws.Range("E26").Value = '  -1.70%  '

$ws.Range("D27").Value = '''6.821'
$ws.Range("E27").Value = '  +5.73%  '

$ws.Range("D28").Value = '''162.56'
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").Value = '''2.478'
$ws.Range("E29").Value = '  -3.14%  '

$ws.Range("D30").Value = '''133.65'

$ws.Range("D31").Value = '''1.124'
$ws.Range("E31").Value = '  -2.93%  '

$ws.Range("D32").Value = '''1.660'
$ws.Range("E32").Value = '  -7.50%  '

$ws.Range("D33").Value = '''0.1047'
$ws.Range("E33").Value = '  -0.73%  '

$ws.Range("D34").Value = '''6.248'
$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").Value = '''3.911'
$ws.Range("E35").Value = '  -1.59%  '

$ws.Range("D36").Value = '''0.02608'
$ws.Range("E36").Value = '  +0.41%  '

$ws.Range("D37").Value = '''9.877'
$ws.Range("E37").Value = '  -9.40%  '

$ws.Range("D38").Value = '''0.06743'
$ws.Range("E38").Value = '  -1.08%  '

$ws.Range("D39").Value = '''12.55'
$ws.Range("E39").Value = '  -1.16%  '

$ws.Range("D40").Value = '''1.341'
$ws.Range("E40").Value = '  -0.66%  '

$ws.Range("D41").Value = '''0.6935'
$ws.Range("E41").Value = '  -1.22%  '

$ws.Range("D42").Value = '''0.2200'
$ws.Range("E42").Value = '  -2.05%  '

$ws.Range("D43").Value = '''0.6732'
$ws.Range("E43").Value = '  -1.89%  '

$ws.Range("D44").Value = '''2.378'
$ws.Range("E44").Value = '  +0.58%  '

$ws.Range("D45").Value = '''14.33'
$ws.Range("E45").Value = '  -0.90%  '

$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").Value = '''1.284'
$ws.Range("E47").Value = '  +5.24%  '

$ws.Range("D48").Value = '''3.634'
$ws.Range("E48").Value = '  -0.23%  '

$ws.Range("D49").Value = '''0.00000000342'
$ws.Range("E49").Value = '  -4.48%  '

$ws.Range("D50").Value = '''1.207'
$ws.Range("E50").Value = '  +2.38%  '

$ws.Range("D51").Value = '''1.208'
$ws.Range("E51").Value = '  -1.36%  '
